$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row above row 12, shifting the MinMaxMed block (old row 12 onward) down.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the MinMaxMed config entry.
$ws.Cells.Item(12, 1).Value = "minMaxMedFilePath"
$ws.Cells.Item(12, 2).Value = "Data\Output\MinMaxMed.xlsx"
$ws.Cells.Item(12, 3).Value = "Location to save the file ""MinMaxMed"""
$ws.Rows.Item(12).RowHeight = 14.25

# Match the final selection left by the editing session.
[void]$ws.Range("C14").Select()
